$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60: small "S6" / "per" label row
$ws.Cells.Item(60, 1).Value = "S6"
$ws.Cells.Item(60, 2).Value = "per"

# Row 61: header row (same labels as rows 12/38/49 groupings)
$ws.Cells.Item(61, 1).Value = "Cores"
$ws.Cells.Item(61, 2).Value = "S_ALMs"
$ws.Cells.Item(61, 3).Value = "S_ALM %"
$ws.Cells.Item(61, 4).Value = "S_Regs"
$ws.Cells.Item(61, 5).Value = "S_Blocks"
$ws.Cells.Item(61, 6).Value = "S_Blocks %"
$ws.Cells.Item(61, 7).Value = "S_Fmax"
$ws.Cells.Item(61, 8).Value = "Fmax Limit"
$ws.Cells.Item(61, 9).Value = "ALM Limit"
$ws.Cells.Item(61, 10).Value = "S_Regs %"

# Row 62
$ws.Cells.Item(62, 1).Value = 1
$ws.Cells.Item(62, 9).Value = 8000
$ws.Range("J62").Formula = "=(D62/16000)*100"

# Row 63 (E63 carries the #,##0.0 number format but stays empty)
$ws.Cells.Item(63, 1).Value = 2
$ws.Range("E63").NumberFormat = "#,##0.0"
$ws.Cells.Item(63, 9).Value = 8000
$ws.Range("J63:J68").Formula = "=(D63/16000)*100"

# Row 64
$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 9).Value = 8000

# Row 65
$ws.Cells.Item(65, 1).Value = 8
$ws.Cells.Item(65, 9).Value = 8000

# Row 66
$ws.Cells.Item(66, 1).Value = 16
$ws.Cells.Item(66, 9).Value = 8000

# Row 67
$ws.Cells.Item(67, 1).Value = 32
$ws.Cells.Item(67, 2).Value = 19294
$ws.Cells.Item(67, 3).Value = 241.18
$ws.Cells.Item(67, 4).Value = 14231
$ws.Cells.Item(67, 5).Value = 10
$ws.Cells.Item(67, 6).Value = 100
$ws.Cells.Item(67, 9).Value = 8000

# Row 68
$ws.Cells.Item(68, 1).Value = 64
$ws.Cells.Item(68, 9).Value = 8000

# Update the active selection to mirror the saved view state
$ws.Range("D23").Select()
